$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row stays with same text (sst table gets reordered by the engine)
$ws.Range("A1").Value = "letter10"
$ws.Range("B1").Value = "letter20"
$ws.Range("C1").Value = "corrAns0"

# Replace "U"/"L" picture-placeholder text with "p1.jpg"/"p2.jpg" image filenames,
# keeping the same U/L assignment pattern (now p1/p2) that lines up with corrAns (c/m)
$p1Rows = @(2,3,6,7,8,11)
$p2Rows = @(4,5,9,10)

foreach ($r in $p1Rows) {
    $ws.Cells.Item($r, 1).Value = "p1.jpg"
    $ws.Cells.Item($r, 2).Value = "p2.jpg"
}
foreach ($r in $p2Rows) {
    $ws.Cells.Item($r, 1).Value = "p2.jpg"
    $ws.Cells.Item($r, 2).Value = "p1.jpg"
}

# Apply a new font (Calibri, 11pt) to the A2:B11 range via a temporary named style,
# then remove the named style so only the underlying cell format (cellXfs) remains.
$tempStyle = $wb.Styles.Add("TempPictureFont")
$tempStyle.Font.Name = "Calibri "
$tempStyle.Font.Size = 11
$ws.Range("A2:B11").Style = "TempPictureFont"
$wb.Styles.Item("TempPictureFont").Delete()

# Selection moved from C11 to B11
$ws.Range("B11").Select()

# Add explicit page setup (portrait, paper size 9 = A4)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
